$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original Text storage type so
# values like "34.50", "559.00", "0.180" keep their exact trailing zeros
# instead of Excel auto-converting the literal into a Number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.683.44"
$ws.Range("E2").Value = "  +5.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.295.99"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.51"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "631.17"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.421"
$ws.Range("E7").Value = "  +10.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.716"
$ws.Range("E8").Value = "  +6.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.295.99"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.591"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000268"
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.180"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.50"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.889.90"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.156.61"
$ws.Range("E17").Value = "  +5.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.275.15"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.25"
$ws.Range("E19").Value = "  +6.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.22"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "438.83"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.01"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.36"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000187"
$ws.Range("E24").Value = "  +44.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.43"
$ws.Range("E25").Value = "  +6.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.17"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.466.08"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.64"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.176"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.73"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "559.00"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.37"
$ws.Range("E34").Value = "  +6.91%  "
$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.73"
$ws.Range("E35").Value = "  +28.04%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.35"
$ws.Range("E36").Value = "  -4.81%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.76"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.36"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.398"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "183.45"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "148.97"
$ws.Range("E46").Value = "  -5.39%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.131"
$ws.Range("E47").Value = "  +9.32%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.99"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.88"
$ws.Range("E50").Value = "  +7.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.19"
$ws.Range("E51").Value = "  -0.97%  "
